$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 244 (existing rows 244:300 shift down to 247:303)
$ws.Rows.Item(244).Resize(3).Insert()

# --- New row 244 ---
$ws.Range("A244").Value = 5
$ws.Range("B244").Value = "Macroferia Regional de Talca"
$ws.Range("C244").Value = "Maule"
$ws.Range("D244").Value = 44798
$ws.Range("E244").Value = 7
$ws.Range("F244").Value = 100112045
$ws.Range("G244").Value = "Zapallo"
$ws.Range("H244").Value = "Camote"
$ws.Range("I244").Value = "1a (guarda)"
$ws.Range("J244").Value = 600
$ws.Range("K244").Value = 750
$ws.Range("L244").Value = 750
$ws.Range("M244").Value = 750
$ws.Range("N244").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O244").Value = "Región del Maule"
$ws.Range("P244").Value = 750
$ws.Range("Q244").Value = 1
$ws.Range("R244").Value = "Hortaliza"

# --- New row 245 ---
$ws.Range("A245").Value = 5
$ws.Range("B245").Value = "Macroferia Regional de Talca"
$ws.Range("C245").Value = "Maule"
$ws.Range("D245").Value = 44798
$ws.Range("E245").Value = 7
$ws.Range("F245").Value = 100112045
$ws.Range("G245").Value = "Zapallo"
$ws.Range("H245").Value = "Camote"
$ws.Range("I245").Value = "2a (guarda)"
$ws.Range("J245").Value = 600
$ws.Range("K245").Value = 600
$ws.Range("L245").Value = 600
$ws.Range("M245").Value = 600
$ws.Range("N245").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O245").Value = "Región del Maule"
$ws.Range("P245").Value = 600
$ws.Range("Q245").Value = 1
$ws.Range("R245").Value = "Hortaliza"

# --- New row 246 ---
$ws.Range("A246").Value = 5
$ws.Range("B246").Value = "Macroferia Regional de Talca"
$ws.Range("C246").Value = "Maule"
$ws.Range("D246").Value = 44798
$ws.Range("E246").Value = 7
$ws.Range("F246").Value = 100112045
$ws.Range("G246").Value = "Zapallo"
$ws.Range("H246").Value = "Paine"
$ws.Range("I246").Value = "1a (guarda)"
$ws.Range("J246").Value = 1500
$ws.Range("K246").Value = 250
$ws.Range("L246").Value = 250
$ws.Range("M246").Value = 250
$ws.Range("N246").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O246").Value = "Región del Maule"
$ws.Range("P246").Value = 250
$ws.Range("Q246").Value = 1
$ws.Range("R246").Value = "Hortaliza"
